# Updated cryptos list on Fri Nov  3 10:51:56 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) and E (Volume) to be treated as plain text so that
# Excel does not coerce values such as "228.82" or "0.0677" into numbers
# (which would change their displayed precision / notation).
$ws.Range("D2:E51").NumberFormat = "@"

$subThree = [string][char]0x2083

# Row 2 - Bitcoin
$ws.Range("D2").Value = "34.496.53"
$ws.Range("E2").Value = "  -2.91%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.797.98"
$ws.Range("E3").Value = "  -2.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.46%  "

# Row 5 - BNB
$ws.Range("D5").Value = "228.82"
$ws.Range("E5").Value = "  -1.20%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  -2.14%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.47%  "

# Row 8 - Solana
$ws.Range("D8").Value = "38.79"
$ws.Range("E8").Value = "  -11.26%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.319"
$ws.Range("E9").Value = "  +2.74%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0675"
$ws.Range("E10").Value = "  -3.88%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -2.17%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.060.68"
$ws.Range("E12").Value = "  -2.26%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "11.07"
$ws.Range("E13").Value = "  -2.17%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.804.42"
$ws.Range("E14").Value = "  -2.09%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.656"
$ws.Range("E15").Value = "  -2.78%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  -4.22%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "34.503.03"
$ws.Range("E17").Value = "  -2.83%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "68.84"
$ws.Range("E18").Value = "  -2.19%  "

# Row 19 - now ShibaInu (was BitcoinCash)
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0" + $subThree + "0777"
$ws.Range("E19").Value = "  -3.28%  "

# Row 20 - now BitcoinCash (was ShibaInu)
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "239.23"
$ws.Range("E20").Value = "  -2.23%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "11.73"
$ws.Range("E21").Value = "  -2.63%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "4.68"
$ws.Range("E22").Value = "  +0.97%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.68%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.55%  "

# Row 25 - Monero
$ws.Range("D25").Value = "172.15"
$ws.Range("E25").Value = "  -0.24%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -4.14%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "17.11"
$ws.Range("E27").Value = "  -4.09%  "

# Row 28 - Stellar
$ws.Range("D28").Value = "0.121"
$ws.Range("E28").Value = "  -1.68%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value = "1.48"
$ws.Range("E29").Value = "  -4.22%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.40%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "4.01"
$ws.Range("E31").Value = "  +1.92%  "

# Row 32 - Hedera
$ws.Range("D32").Value = "0.0539"
$ws.Range("E32").Value = "  -2.08%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "3.87"
$ws.Range("E33").Value = "  -5.38%  "

# Row 34 - TrustWalletToken
$ws.Range("D34").Value = "1.23"
$ws.Range("E34").Value = "  +8.27%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "1.78"
$ws.Range("E35").Value = "  -3.32%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.689"
$ws.Range("E36").Value = "  -0.52%  "

# Row 37 - Aave
$ws.Range("D37").Value = "90.67"
$ws.Range("E37").Value = "  -4.89%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  +4.21%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.312.26"
$ws.Range("E39").Value = "  -2.71%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -2.68%  "

# Row 41 - now ARBITRUM (was HuobiToken)
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "0.953"
$ws.Range("E41").Value = "  -5.82%  "

# Row 42 - now HuobiToken (was ARBITRUM)
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "2.43"
$ws.Range("E42").Value = "  -1.41%  "

# Row 43 - InjectiveProtocol
$ws.Range("D43").Value = "14.16"
$ws.Range("E43").Value = "  -8.16%  "

# Row 44 - RenderToken
$ws.Range("D44").Value = "2.20"
$ws.Range("E44").Value = "  -10.58%  "

# Row 45 - MXToken
$ws.Range("E45").Value = "  -4.26%  "

# Row 46 - FraxShare
$ws.Range("D46").Value = "6.19"
$ws.Range("E46").Value = "  -1.34%  "

# Row 47 - Kaspa
$ws.Range("E47").Value = "  -1.08%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "1.982.59"
$ws.Range("E48").Value = "  -1.50%  "

# Row 49 - PaxDollar
$ws.Range("E49").Value = "  +0.48%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.0659"
$ws.Range("E50").Value = "  +3.67%  "

# Row 51 - Quant
$ws.Range("D51").Value = "97.22"
$ws.Range("E51").Value = "  -5.01%  "
